$d = $word.ActiveDocument

# 1) "e ffisselle de " -> "e fisselle de " (fix double-f typo)
$d.Content.Find.Execute("e ffisselle de ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e fisselle de ", 2)

# 2) "moing" -> "moingdre" (extend word)
$d.Content.Find.Execute("moing", $true, $false, $false, $false, $false,
                         $true, 1, $false, "moingdre", 2)

# 3) "<ill/></ms>" -> "</ms>" (drop <ill/> tag)
$d.Content.Find.Execute("<ill/></ms>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "</ms>", 2)
